$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the database table (rows 16-21) with the new account-statement data.
# Row 16 (LINO - 1805) stays as-is.
# The remaining rows are reordered/refreshed per worker with the updated
# "Salario Basico" (column G) figures for CARLOS and JORGE.

# Row 17: CARLOS ANDRES LEON FRANCO - periodo 1804
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047457155"
$ws.Range("D17").Value = "CARLOS ANDRES LEON FRANCO"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 737717

# Row 18: JORGE RAFAEL LEON FRANCO - periodo 1804
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047471570"
$ws.Range("D18").Value = "JORGE RAFAEL LEON FRANCO"
$ws.Range("E18").Value = "1804"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 737717

# Row 19: LINO RICARDO LEON BOLIVAR - periodo 1805
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73145513"
$ws.Range("D19").Value = "LINO RICARDO LEON BOLIVAR"
$ws.Range("E19").Value = "1805"
$ws.Range("F19").Value = 29269
$ws.Range("G19").Value = 731717

# Row 20: CARLOS ANDRES LEON FRANCO - periodo 1805
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047457155"
$ws.Range("D20").Value = "CARLOS ANDRES LEON FRANCO"
$ws.Range("E20").Value = "1805"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 737717

# Row 21: JORGE RAFAEL LEON FRANCO - periodo 1805
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047471570"
$ws.Range("D21").Value = "JORGE RAFAEL LEON FRANCO"
$ws.Range("E21").Value = "1805"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 737717
